# #5: cash & deposit done
# Add bank / deposit_type / currency (+ shared metadata columns) to the
# "存款" (deposit) sheet, matching the schema used on the other sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1) -----------------------------------------------
# Existing B1:F1 held stale literal copies of row-2 values; relabel them
# as proper column headers and extend with the shared metadata headers.
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# --- Data rows (rows 2-7) ----------------------------------------------
# F2 used to be the text "1300849" (shared string); make it a real number.
$ws.Range("F2").Value = 1300849

# Column I holds the text date "2013-12-26" - force a text number format
# first so Excel doesn't silently reinterpret the literal as a date serial.
$ws.Range("I2:I7").NumberFormat = "@"

$rows = 2,3,4,5,6,7
$indexes = 49,50,51,52,53,54
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $idx = $indexes[$i]
    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"
    $ws.Range("I$r").Value = "2013-12-26"
    $ws.Range("J$r").Value = "潘孟安"
    $ws.Range("K$r").Value = 1376
    $ws.Range("L$r").Value = "tmpf07c1"
    $ws.Range("M$r").Value = $idx
}

# --- Formatting ----------------------------------------------------------
# New header cells G1:M1 need the same bold/bordered look as B1:F1; new
# data cells G2:M7 need the same plain look as the rest of the data rows.
# Copy formats last so the text-vs-date coercion above already happened.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122)

$ws.Range("B2").Copy() | Out-Null
$ws.Range("G2:M7").PasteSpecial(-4122)
